$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 2 4 "26.202.33"
Set-TextValue 2 5 "  +1.54%  "
Set-TextValue 3 4 "1.605.10"
Set-TextValue 3 5 "  +0.34%  "
Set-TextValue 4 5 "  -0.27%  "
Set-TextValue 5 4 "212.23"
Set-TextValue 5 5 "  +1.60%  "
Set-TextValue 6 5 "  -0.24%  "
Set-TextValue 7 4 "0.482"
Set-TextValue 7 5 "  +0.41%  "
Set-TextValue 8 5 "  +1.14%  "
Set-TextValue 9 5 "  +1.25%  "
Set-TextValue 10 4 "18.16"
Set-TextValue 10 5 "  +1.32%  "
Set-TextValue 11 4 "0.0797"
Set-TextValue 11 5 "  +1.64%  "
Set-TextValue 12 4 "1.826.72"
Set-TextValue 12 5 "  +0.23%  "
Set-TextValue 13 4 "1.606.92"
Set-TextValue 13 5 "  +0.30%  "
Set-TextValue 14 4 "4.01"
Set-TextValue 14 5 "  -0.99%  "
Set-TextValue 15 4 "0.509"
Set-TextValue 15 5 "  -0.14%  "
Set-TextValue 16 4 "26.174.27"
Set-TextValue 16 5 "  +1.41%  "
Set-TextValue 17 4 "60.64"
Set-TextValue 17 5 "  +0.22%  "
Set-TextValue 18 4 "0.0₃0728"
Set-TextValue 18 5 "  +1.61%  "
Set-TextValue 19 5 "  -0.12%  "
Set-TextValue 20 4 "198.82"
Set-TextValue 20 5 "  +4.88%  "
Set-TextValue 21 4 "4.24"
Set-TextValue 21 5 "  +1.26%  "
Set-TextValue 22 4 "9.39"
Set-TextValue 22 5 "  +0.53%  "
Set-TextValue 23 5 "  +1.08%  "
Set-TextValue 24 4 "0.132"
Set-TextValue 24 5 "  +2.24%  "
Set-TextValue 25 4 "142.13"
Set-TextValue 25 5 "  +0.98%  "
Set-TextValue 26 5 "  +3.01%  "
Set-TextValue 27 5 "  -0.30%  "
Set-TextValue 28 4 "15.18"
Set-TextValue 28 5 "  +1.33%  "
Set-TextValue 29 5 "  -0.46%  "
Set-TextValue 30 5 "  -1.06%  "
Set-TextValue 31 5 "  +0.79%  "
Set-TextValue 32 5 "  +1.81%  "
Set-TextValue 33 5 "  +0.54%  "
Set-TextValue 34 5 "  +1.86%  "
Set-TextValue 35 5 "  -1.55%  "
Set-TextValue 36 4 "1.107.83"
Set-TextValue 36 5 "  +1.12%  "
Set-TextValue 37 2 "MXToken"
Set-TextValue 37 3 "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue 37 4 "2.35"
Set-TextValue 37 5 "  -0.53%  "
Set-TextValue 38 2 "PaxDollar"
Set-TextValue 38 3 "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue 38 4 "1.00"
Set-TextValue 38 5 "  +0.12%  "
Set-TextValue 39 5 "  +0.67%  "
Set-TextValue 40 4 "0.502"
Set-TextValue 40 5 "  +0.82%  "
Set-TextValue 41 4 "0.786"
Set-TextValue 41 5 "  -0.77%  "
Set-TextValue 42 4 "0.778"
Set-TextValue 42 5 "  +4.66%  "
Set-TextValue 43 4 "1.739.23"
Set-TextValue 43 5 "  +0.21%  "
Set-TextValue 44 5 "  +0.99%  "
Set-TextValue 45 4 "92.74"
Set-TextValue 45 5 "  -2.98%  "
Set-TextValue 46 4 "1.55"
Set-TextValue 46 5 "  +8.73%  "
Set-TextValue 47 5 "  -7.65%  "
Set-TextValue 48 4 "53.53"
Set-TextValue 48 5 "  +0.52%  "
Set-TextValue 49 5 "  -0.46%  "
Set-TextValue 50 5 "  -0.14%  "
Set-TextValue 51 5 "  +0.05%  "
